$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date value "02/08/2023" -> "05/08/2023" wherever it appears,
# preserving the text "quote prefix" formatting already on these cells.
foreach ($cell in $ws.Range("B1:B11").Cells) {
    if ($cell.Text -eq "02/08/2023") {
        $cell.Value = "'05/08/2023"
    }
}

# Change the active cell selection to H9
$ws.Range("H9").Select()

# Configure page setup (paper size 9, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
